$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.843037247657776
$ws.Range("B1").Value = 2.067568063735962
$ws.Range("C1").Value = 2.25141978263855
$ws.Range("D1").Value = 3.283270120620728
$ws.Range("E1").Value = 1.769476294517517
